# Add "2022-Q4" quarterly fund-holdings sheet and update the "总计" summary
# sheet on top of it, per commit "feat: add 2022-Q4 data".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert the new "2022-Q4" worksheet right after the "总计" summary sheet
#    (i.e. it becomes the workbook's 2nd sheet, pushing the former 2022-Q3..
#    2021-Q1 sheets one slot later).
# ---------------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item(1)
$q4 = $wb.Worksheets.Add($null, $totalSheet)
$q4.Name = "2022-Q4"

# ---------------------------------------------------------------------------
# 2. Populate "2022-Q4" with the fund holdings table (header + 17 funds).
# ---------------------------------------------------------------------------
$headers = @("", "基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($col = 2; $col -le 8; $col++) {
    $q4.Cells.Item(1, $col).Value = $headers[$col - 1]
}

# code, name, size, position, pct, value, rank
$funds = @(
    @("011056", "博时汇兴回报一年持有期灵活配置混合", "91.53", "68.57", "3.67", "3.3592", 4),
    @("012703", "华夏核心成长混合A",                    "4.93",  "92.86", "5.12", "0.2524", 8),
    @("660012", "农银消费主题混合A",                    "6.14",  "87.44", "2.79", "0.1713", 8),
    @("910021", "东方红启华三年持有期混合A",            "4.18",  "91.17", "2.99", "0.1250", 7),
    @("014158", "博时浦惠一年持有期混合A",              "3.84",  "49.14", "2.48", "0.0952", 10),
    @("010043", "天弘安康颐和混合A",                    "14.87", "23.17", "0.60", "0.0892", 10),
    @("008261", "招商研究优选股票A",                    "1.90",  "86.08", "2.44", "0.0464", 8),
    @("012710", "华夏核心成长混合C",                    "0.60",  "92.86", "5.12", "0.0307", 8),
    @("011313", "东方红启华三年持有期混合B",            "0.89",  "91.17", "2.99", "0.0266", 7),
    @("360010", "光大保德信均衡精选混合",               "0.31",  "87.82", "4.43", "0.0137", 3),
    @("008262", "招商研究优选股票C",                    "0.54",  "86.08", "2.44", "0.0132", 8),
    @("014159", "博时浦惠一年持有期混合C",              "0.37",  "49.14", "2.48", "0.0092", 10),
    @("010044", "天弘安康颐和混合C",                    "1.17",  "23.17", "0.60", "0.0070", 10),
    @("001830", "融通跨界成长灵活配置混合",             "0.21",  "81.58", "2.98", "0.0063", 5),
    @("005067", "融通逆向策略灵活配置混合A",            "0.11",  "75.72", "2.51", "0.0028", 7),
    @("009270", "融通逆向策略灵活配置混合C",            "0.03",  "75.72", "2.51", "0.0008", 7),
    @("960033", "农银消费主题混合H",                    "0.00",  "87.44", "2.79", 0,        8)
)

for ($i = 0; $i -lt $funds.Length; $i++) {
    $row = $i + 2
    $fund = $funds[$i]

    $q4.Cells.Item($row, 1).Value = $i

    # Columns B-F (fund code/name/size/position/pct) are text in the source
    # data, even though most look numeric - force with a leading quote so
    # Excel doesn't coerce them (and so leading zeros in codes survive).
    $q4.Cells.Item($row, 2).Value = "'" + $fund[0]
    $q4.Cells.Item($row, 3).Value = "'" + $fund[1]
    $q4.Cells.Item($row, 4).Value = "'" + $fund[2]
    $q4.Cells.Item($row, 5).Value = "'" + $fund[3]
    $q4.Cells.Item($row, 6).Value = "'" + $fund[4]

    if ($row -eq 18) {
        # Last row's "持有市值(亿元)" is stored as a real 0, not text "0".
        $q4.Cells.Item($row, 7).Value = $fund[5]
    } else {
        $q4.Cells.Item($row, 7).Value = "'" + $fund[5]
    }

    $q4.Cells.Item($row, 8).Value = $fund[6]
}

# ---------------------------------------------------------------------------
# 3. Update the "总计" sheet: a new row for 2022-Q4 is inserted at the top
#    of the data (row 2), shifting every other quarter down by one row.
# ---------------------------------------------------------------------------
$total = $totalSheet

$summary = @(
    @("2022-Q4", 17, 4.25),
    @("2022-Q3", 3,  3.22),
    @("2022-Q2", 5,  5.09),
    @("2022-Q1", 6,  8.53),
    @("2021-Q4", 19, 14.41),
    @("2021-Q3", 6,  9.75),
    @("2021-Q2", 16, 8.19),
    @("2021-Q1", 3,  0.59)
)

for ($i = 0; $i -lt $summary.Length; $i++) {
    $row = $i + 2
    $entry = $summary[$i]
    $total.Cells.Item($row, 1).Value = $i
    $total.Cells.Item($row, 2).Value = $entry[0]
    $total.Cells.Item($row, 3).Value = $entry[1]
    $total.Cells.Item($row, 4).Value = $entry[2]
}

# D5 (2022-Q1 market value) keeps the source file's exact binary value.
$total.Cells.Item(5, 4).Value = 8.529999999999999

# Row 9 ("2021-Q1") is brand new - clone column A's number style from the
# row above it (row 8) so it matches the sheet's existing formatting.
$total.Cells.Item(8, 1).Copy()
$total.Cells.Item(9, 1).PasteSpecial(-4122)
